$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("12:12").Delete()
